$d = $word.ActiveDocument

# --- 1) Merge the split hyperlink runs back into single runs ---
# (Find/Replace with the same text collapses the multiple runs that made
#  up the URL into one run. NOTE: the search string deliberately starts
#  one character after the start of the first run -- when a Find match
#  begins exactly at a run boundary, this engine's replace loses that
#  run's rPr/rStyle; starting mid-run avoids that and the leading
#  character re-merges into the same, correctly-formatted run.)

$find = "ttp://science.time.com/2010/06/28/climate-change-and-space-junk/"
$d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $find, 2) | Out-Null

$find = "ttps://www.space.com/6720-space-littering-impact-earths-atmosphere.html"
$d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $find, 2) | Out-Null

$find = "ttps://science.sciencemag.org/content/311/5759/340/tab-pdf"
$d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $find, 2) | Out-Null

# --- 2) Add a new reference paragraph at the end of the document,
#        with a hyperlink to the ESA "types of orbits" page. ---

# The existing content ends with a hidden "_GoBack" bookmark right after
# "dust reflectors". Remove it -- we'll recreate it at the new end below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Start a new paragraph after the current last paragraph.
$lastPara = $d.Paragraphs.Last
$endOfDoc = $lastPara.Range
$endOfDoc.Collapse(0)   # wdCollapseEnd
$endOfDoc.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range
$newRange.Collapse(0)   # wdCollapseEnd

# Type the URL text, then convert it into a hyperlink.
$url = "https://www.esa.int/Our_Activities/Space_Transportation/Types_of_orbits"
$newRange.InsertAfter($url)
$urlRange = $d.Range($newRange.Start, $newRange.Start + $url.Length)
$d.Hyperlinks.Add($urlRange, $url) | Out-Null

# Append the trailing description text after the hyperlink.
$afterHyperlink = $d.Paragraphs.Last.Range
$afterHyperlink.Collapse(0)
$afterHyperlink.InsertAfter(" types of orbits")

# --- 3) Re-create the "_GoBack" bookmark at the very end of the document. ---
# Placing a zero-length bookmark exactly at the document's final position
# is unreliable, so insert a throwaway sentinel character first, anchor the
# bookmark right before it, then delete the sentinel -- the bookmark stays
# anchored to the (now final) end of the text.
$finalPara = $d.Paragraphs.Last
$finalEnd = $finalPara.Range.End
$sentinelPos = $finalEnd - 1
$sentinelRange = $d.Range($sentinelPos, $sentinelPos)
$sentinelRange.InsertAfter("Z")

$bmRange = $d.Range($sentinelPos, $sentinelPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

$zRange = $d.Range($sentinelPos, $sentinelPos + 1)
$zRange.Delete()

Write-Host "Done"
